# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.415.06'
$ws.Range("E2").Value = '  +1.21%  '

# Row 3
$ws.Range("D3").Value = '3.018.03'
$ws.Range("E3").Value = '  +2.31%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.54'
$ws.Range("E5").Value = '  +2.14%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.40'
$ws.Range("E6").Value = '  +4.96%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = '3.016.92'
$ws.Range("E8").Value = '  +2.41%  '

# Row 9
$ws.Range("E9").Value = '  -0.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.38'
$ws.Range("E10").Value = '  +11.77%  '

# Row 11
$ws.Range("E11").Value = '  +3.98%  '

# Row 12
$ws.Range("E12").Value = '  +0.28%  '

# Row 13
$ws.Range("E13").Value = '  +3.61%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.46'
$ws.Range("E14").Value = '  +1.41%  '

# Row 15
$ws.Range("E15").Value = '  +2.69%  '

# Row 16
$ws.Range("D16").Value = '3.519.04'
$ws.Range("E16").Value = '  +2.15%  '

# Row 17
$ws.Range("E17").Value = '  +0.69%  '

# Row 18
$ws.Range("D18").Value = '62.366.70'
$ws.Range("E18").Value = '  +1.10%  '

# Row 19
$ws.Range("D19").Value = '3.020.54'
$ws.Range("E19").Value = '  +2.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.43'
$ws.Range("E20").Value = '  -0.05%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.21'
$ws.Range("E21").Value = '  +2.97%  '

# Row 22
$ws.Range("E22").Value = '  +1.90%  '

# Row 23
$ws.Range("E23").Value = '  +2.51%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.30'
$ws.Range("E24").Value = '  +1.54%  '

# Row 25
$ws.Range("E25").Value = '  +4.49%  '

# Row 26
$ws.Range("E26").Value = '  +13.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("E27").Value = '  -0.60%  '

# Row 28
$ws.Range("E28").Value = '  +0.03%  '

# Row 29
$ws.Range("E29").Value = '  +3.80%  '

# Row 30
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  +4.51%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.14'
$ws.Range("E32").Value = '  +4.13%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.58'

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  +2.84%  '

# Row 35
$ws.Range("D35").Value = '0.0₃0847'
$ws.Range("E35").Value = '  +9.44%  '

# Row 36
$ws.Range("E36").Value = '  +2.65%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.85'
$ws.Range("E37").Value = '  +3.49%  '

# Row 38
$ws.Range("E38").Value = '  -0.01%  '

# Row 39
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.18'
$ws.Range("E39").Value = '  +0.40%  '

# Row 40
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").Value = '  +7.52%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.03'
$ws.Range("E41").Value = '  -0.49%  '

# Row 42
$ws.Range("E42").Value = '  +1.87%  '

# Row 43
$ws.Range("E43").Value = '  +8.26%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.32'
$ws.Range("E44").Value = '  +10.32%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '391.20'
$ws.Range("E45").Value = '  +0.81%  '

# Row 46
$ws.Range("E46").Value = '  +0.63%  '

# Row 47
$ws.Range("D47").Value = '2.755.11'
$ws.Range("E47").Value = '  +1.59%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.58'
$ws.Range("E48").Value = '  +2.92%  '

# Row 49
$ws.Range("E49").Value = '  +0.08%  '

# Row 50
$ws.Range("E50").Value = '  +1.36%  '

# Row 51
$ws.Range("E51").Value = '  +0.04%  '
